# Atualizacao rapida de agenda as 9:45:49,74
# Reorganizes the Giovani/Roberto technician agenda rows (2-11) with
# fresh statuses/observations and updates the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: column B (ID) cells use a "quote-prefixed text" number format in
# this sheet (style index 9, e.g. '0891) so a plain digit-string value
# renders like a number. Writing with a leading apostrophe keeps Excel's
# text-entry semantics and preserves that style (COM strips the
# apostrophe from the stored value automatically).

# --- Row 2 (Giovani / Acid) ---
$ws.Range("A2").Value = "Giovani"
$ws.Range("B2").Value = "'0891"
$ws.Range("C2").Value = "Acid"
$ws.Range("D2").Value = "Câmera 2 parece que caiu, passar para o DDNS também."
$ws.Range("G2").Value = "Pendente"
$ws.Range("H2").Value = "Maxvel: 22 / Forte: 12"

# --- Row 3 (Giovani / MegaScan) ---
$ws.Range("A3").Value = "Giovani"
$ws.Range("B3").Value = "'0643"
$ws.Range("C3").Value = "MegaScan"
$ws.Range("D3").Value = "Sem comunicação de câmeras, é via DDNS."
$ws.Range("G3").Value = "Pendente"

# --- Row 4 (Giovani / Escola Manoel Correira) ---
$ws.Range("A4").Value = "Giovani"
$ws.Range("B4").Value = "'0756"
$ws.Range("C4").Value = "Escola Manoel Correira"
$ws.Range("D4").Value = "Sem comunicação de câmeras, favor restaurar e passar para o DDNS."
$ws.Range("G4").Value = "Pendente"

# --- Row 5 (Giovani / Sitio Alves) ---
$ws.Range("A5").Value = "Giovani"
$ws.Range("B5").Value = "'0422"
$ws.Range("C5").Value = "Sitio Alves"
$ws.Range("D5").Value = "Sem comunicação de câmeras, central tá no gprs e zona aberta. Cliente pedindo reparo em tudo e favor passar para o DDNS."
$ws.Range("G5").Value = "Pendente"

# --- Row 6 (Giovani / Cetep) ---
$ws.Range("A6").Value = "Giovani"
$ws.Range("B6").Value = "'0790"
$ws.Range("C6").Value = "Cetep"
$ws.Range("D6").Value = "Câmeras muito instáveis, cliente pedindo que arrume. Favor passar para o DDNS."
$ws.Range("G6").Value = "Pendente"

# --- Row 7 (Roberto / Galpão Toyota) ---
$ws.Range("A7").Value = "Roberto"
$ws.Range("B7").Value = "'0803"
$ws.Range("C7").Value = "Galpão Toyota"
$ws.Range("D7").Value = "Parece que o local agora tem internet, colocar central via internet. Se tiver câmera, programar pra gente monitrar. Se não tiver, instalar e programar. Local com muitos disparos frequentes também."
$ws.Range("D7").WrapText = $True
$ws.Range("G7").Value = "Pendente"

# --- Row 8 (Roberto / BC Refratário) ---
$ws.Range("A8").Value = "Roberto"
$ws.Range("B8").Value = "'0463"
$ws.Range("C8").Value = "BC Refratário"
$ws.Range("D8").Value = "Zona aberta, colocar câmeras via DDNS (antes era)."
$ws.Range("G8").Value = "Pendente"

# --- Row 9 (Roberto / Recapel) ---
$ws.Range("A9").Value = "Roberto"
$ws.Range("B9").Value = "'0869"
$ws.Range("C9").Value = "Recapel"
$ws.Range("D9").Value = "Registrado que o sensor quebrou, favor arrumar. Limpeza na câmera 8 e passar as câmeras para o DDNS."
$ws.Range("D9").WrapText = $True
$ws.Range("G9").Value = "Pendente"

# --- Row 10 (Roberto / Med Center) ---
$ws.Range("A10").Value = "Roberto"
$ws.Range("B10").Value = "'0079"
$ws.Range("C10").Value = "Med Center"
$ws.Range("D10").Value = "Sem comunicação de alarmes, funciona via internet."
$ws.Range("G10").Value = "Pendente"

# --- Row 11 (Roberto / Brapi) -- was a blank row, now filled in ---
$ws.Range("A11").Value = "Roberto"
$ws.Range("B11").Value = "'0217"
$ws.Range("C11").Value = "Brapi"
$ws.Range("D11").Value = "Zonas abertas, é AMT 8000. Pedro diretor pediu pra arrumar."
$ws.Range("G11").Value = "Pendente"

# --- Row heights: content reflow changed which observations wrap onto
#     2 or 3 lines at the current column width, so several rows' auto
#     heights move around. Rows going back to the single-line default
#     are auto-fitted; rows that now need extra wrapped lines get an
#     explicit height. ---
$ws.Rows("3:3").AutoFit()
$ws.Rows("4:4").AutoFit()
$ws.Rows("5:5").RowHeight = 30
$ws.Rows("7:7").RowHeight = 45
$ws.Rows("8:8").AutoFit()
$ws.Rows("9:9").RowHeight = 30

# --- Selection moved to H11 (last touched cell) ---
$ws.Range("H11").Select()
